$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Log")

# ---------------------------------------------------------------------------
# Row 16 (Id 15): status flips from "Open" to "Closed" and the yellow
# highlight used while the row was open is cleared.
# ---------------------------------------------------------------------------
$row16 = $ws.Range("A16:I16")
$row16.Interior.ColorIndex = -4142
$row16.Interior.Pattern = -4142
$ws.Cells.Item(16, 9).Value = "Closed"

# ---------------------------------------------------------------------------
# Row 26 (Id 25): same treatment - status flips to "Closed" and the row
# loses its highlight now that the task is done.
# ---------------------------------------------------------------------------
$row26 = $ws.Range("A26:I26")
$row26.Interior.ColorIndex = -4142
$row26.Interior.Pattern = -4142
$ws.Cells.Item(26, 9).Value = "Closed"

# ---------------------------------------------------------------------------
# Row 27: brand-new log entry recording the creation of the
# "MinutesOfMeeting - w15.3" document.
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 41014
$ws.Cells.Item(27, 3).Value = "Action"
$ws.Cells.Item(27, 4).Value = 'Create "MinutesOfMeeting - w15.3" document'
$ws.Cells.Item(27, 5).Value = "Mikael"
$ws.Cells.Item(27, 6).Value = 41015
$ws.Cells.Item(27, 8).Value = "Yes"
$ws.Cells.Item(27, 9).Value = "Closed"

# ---------------------------------------------------------------------------
# Scroll / selection follows the newly added row.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("A27").Select()
